# "The Last Update 15-03-2024"
# Refresh the Bundesliga standings table: team order shifts slightly (rows
# 7-15) and every stat column (Cartões, Escanteios, 1.5+, 2.5+, Med. Gols)
# gets new values. The "#" rank column (A) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Equipe, Cartões, Escanteios, "1.5+", "2.5+", "Med. Gols")
$data = @{
    2  = @("Leverkusen",    "1.8", "6.4", "92%", "75%", "3.12")
    3  = @("Bayern",        "1.7", "7.0", "88%", "88%", "4.00")
    4  = @("Stuttgart",     "1.6", "4.7", "96%", "75%", "3.56")
    5  = @("Dortmund",      "1.7", "4.9", "83%", "58%", "3.20")
    6  = @("RB Leipzig",    "1.8", "5.9", "88%", "71%", "3.50")
    7  = @("Frankfurt",     "1.9", "4.0", "83%", "53%", "2.84")
    8  = @("Hoffenheim",    "2.7", "4.2", "96%", "75%", "3.52")
    9  = @("Werder Bremen", "2.0", "3.1", "79%", "54%", "2.96")
    10 = @("Freiburg",      "1.4", "3.4", "79%", "67%", "3.14")
    11 = @("Augsburg",      "2.4", "4.1", "92%", "63%", "3.40")
    12 = @("Heidenheim",    "1.7", "5.6", "83%", "58%", "3.17")
    13 = @("M'Gladbach",    "2.0", "5.3", "88%", "79%", "3.67")
    14 = @("Wolfsburg",     "2.3", "3.8", "83%", "58%", "2.78")
    15 = @("Union Berlin",  "1.8", "4.4", "75%", "46%", "2.56")
    16 = @("Bochum",        "3.0", "3.5", "88%", "63%", "3.24")
    17 = @("FC Köln",       "1.7", "4.9", "79%", "29%", "2.23")
    18 = @("Mainz",         "2.8", "5.2", "71%", "42%", "2.32")
    19 = @("Darmstadt",     "2.5", "3.7", "75%", "62%", "3.42")
}

# Columns C:G hold numeric-looking text ("1.8", "92%", "3.12", ...). Excel's
# automatic type inference would otherwise silently convert those into real
# numbers/percentages, which changes the cell type away from the shared
# string it needs to stay as. Mark the range as Text first, write the
# values, then drop back to the Normal style so no stray number format
# lingers on the cells (matches the original workbook, which carried no
# explicit style on these data cells).
$statRange = $ws.Range("C2:G19")
$statRange.NumberFormat = "@"

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

$statRange.Style = "Normal"
